$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supercoiling row in the "Ancestrally Essential" table (row 23):
# Previously split across C23 ("gyrAB, parE") and D23 ("parC"); the new
# gene clustering merges them into a single value matching the Core
# Essential table's Supercoiling entry, and D23 is cleared.
$ws.Range("C23").Value = "gyrAB, parCE"
$ws.Range("D23").Value = ""

# "Other" row in the "Ancestrally Essential" table (row 33): rodZ was
# added to the ancestrally essential gene clustering and nrfF was removed.
$ws.Range("C33").Value = "acpPS, adk, asd, bamAD, birA, cca, cdsA, coaADE, cohE, csrA, dapABDE, def, der, dfp, dnaAK, dut, dxr, dxs, eno, era, erpA, fbaA, ffh, fldA, fmt, folABCDEK, frr, ftsBI, gapA, glmMSU, glyA, gmk, gpsA, groLS, grpE, hemABCDEGHL, iscS, ispABDEFGHU, kdsAB, lepB, lexA, lgt, ligA, lnt, lolABCDE, lpd, lptABDEFG, lpxABCDHK, lspA, metK, mnmA, mraY, mrdAB, mreBCD, msbA, mukEF, murJ, nadDEK, nrdAB, obgE, orn, pgk, pgsA, plsBC, ppa, prs, psd, pssA, pth, pyrGH, ribABCDEF, rimM, rnc, rnpA, rodZ, secADEFY, suhB, tadA, thiL, thyA, tilS, tmk, topA, trmD, tsaBDE, ubiABDEGHX, waaA, ybeY, ygfZ, yidC, yihA, yqgF"
